$wb = $excel.ActiveWorkbook

$setup = $wb.Worksheets.Item("Setup")
$createFieldData = $wb.Worksheets.Item("createFieldData")

# Update Setup sheet data: row 6 gets a new DATA field entry
# (Table=1, FieldNo=35, Just=l), and row 5 Field changes from DATA to ID.
# Order matters for shared string table indices: "l" must be registered
# before "ID".
$setup.Range("B6").Value = "DATA"
$setup.Range("C6").Value = 1
$setup.Range("D6").Value = 35
$setup.Range("E6").Value = "l"

$setup.Range("B5").Value = "ID"

# Update selections on each sheet
$setup.Range("B6").Select()
$createFieldData.Range("C4:C5").Select()
$createFieldData.Range("C5").Activate()

# Make createFieldData the active sheet/tab
$createFieldData.Activate()

$wb.Save()
